# Add team record (Wins / Losses / Ties) columns to the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1) - new columns AD, AE, AF
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Match the bold/centered/bordered header style used by the other header
# cells (e.g. AC1) by copying its formatting across.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# Data rows 2-47 all get the same team record: 83 wins, 79 losses, 0 ties
for ($r = 2; $r -le 47; $r++) {
    $ws.Cells.Item($r, 30).Value = 83   # AD
    $ws.Cells.Item($r, 31).Value = 79   # AE
    $ws.Cells.Item($r, 32).Value = 0    # AF
}
